$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted before the existing row 392,
# shifting all subsequent rows (392-480) down by one (to 393-481).
$ws.Rows("392:392").Insert()

# Populate the newly inserted row 392 with the new record's data.
$ws.Range("A392").Value = 8
$ws.Range("B392").Value = "Terminal La Palmera de La Serena"
$ws.Range("C392").Value = "Coquimbo"
$ws.Range("D392").Value = 45258
$ws.Range("E392").Value = 4
$ws.Range("F392").Value = 100112031
$ws.Range("G392").Value = "Poroto verde"
$ws.Range("H392").Value = "Sin especificar"
$ws.Range("I392").Value = "Primera"
$ws.Range("J392").Value = 400
$ws.Range("K392").Value = 41000
$ws.Range("L392").Value = 42000
$ws.Range("M392").Value = 41500
$ws.Range("N392").Value = "`$/malla 25 kilos"
$ws.Range("O392").Value = "Región de Arica y Parinacota"
$ws.Range("P392").Value = 1660
$ws.Range("Q392").Value = 25
$ws.Range("R392").Value = "Hortaliza"
